$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I0 and IF, matching style of existing headers
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$I0 = @(8, 9, 4, 6, 4, 6, 5, 5, 8, 8, 1, 7)
$IF = @(8, 9, 5, 6, 4, 6, 5, 6, 9, 8, 4, 7)

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I0[$i]
    $ws.Cells.Item($row, 10).Value = $IF[$i]
}
